$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cr = [char]13
$ws.Cells.Item(1, 1).Value = 'Closures/remote learning and relocationsGoogle Tag Manager (noscript)End Google Tag Manager (noscript)You may be trying to access this site from a secured browser on the server. Please enable scripts and reload this page.Start of Site headerSkip to contentSearchSearchEnd of Site headerHome>The Department>Program directory>Emergencies and Natural Disasters>Closures/remote learning and relocationsClosures/remote learning and relocations- Example page content area using current sharepoint wrapperPage ContentThis page lists early childhood services, schools and TAFEs that are currently closed.SeeCoronavirus advicefor the latest advice on the outbreak of COVID-19 (coronavirus).On this pageSchool and early childhood service, TAFE closures and relocationsBus service cancellations or alterations*Please note these closures are NOT all related to the Covid-19 pandemic.School and early childhood service, TAFE closures and relocations for Thursday 23 July 2020South-Eastern Victoria RegionEarly childhood servicesThe Department has been advised of the followingearly childhood serviceclosures:'
$ws.Cells.Item(2, 1).Value = 'li: Ada Mary A''beckett Children''s Centre Inc PORT MELBOURNE'
$ws.Cells.Item(3, 1).Value = 'li: Aqua Energy Creche SALE'
$ws.Cells.Item(4, 1).Value = 'li: Berwick Neighbourhood Centre - Marriott Waters LYNDHURST'
$ws.Cells.Item(5, 1).Value = 'li: Berwick Neighbourhood Centre (Timbarra Playroom) BERWICK'
$ws.Cells.Item(6, 1).Value = 'li: Brighton Grammar OSHClub BRIGHTON'
$ws.Cells.Item(7, 1).Value = 'li: Camp Australia - Geelong Grammar School Toorak Campus OSHC TOORAK'
$ws.Cells.Item(8, 1).Value = 'li: Camp Australia - Haileybury Edrington OSHC BERWICK'
$ws.Cells.Item(9, 1).Value = 'li: Camp Australia - Haileybury Newlands OSHC KEYSBOROUGH'
$ws.Cells.Item(10, 1).Value = 'li: Camp Australia - Loreto College Mandeville Hall OSHC TOORAK'
$ws.Cells.Item(11, 1).Value = 'li: Camp Australia - Mentone Grammar OSHC MENTONE'
$ws.Cells.Item(12, 1).Value = 'li: Camp Australia - Peninsula Grammar OSHC MOUNT ELIZA'
$ws.Cells.Item(13, 1).Value = 'li: Camp Australia - Sacred Heart Parish School - Sandringham OSHC  SANDRINGHAM'
$ws.Cells.Item(14, 1).Value = 'li: Camp Australia - St Kevins College Toorak OSHC TOORAK'
$ws.Cells.Item(15, 1).Value = 'li: Camp Australia - Toorak College OSHC MOUNT ELIZA'
$ws.Cells.Item(16, 1).Value = 'li: Carnegie Occasional Care Centre CARNEGIE'
$ws.Cells.Item(17, 1).Value = 'li: Castlefield Community Centre HAMPTON'
$ws.Cells.Item(18, 1).Value = 'li: Cranbourne Community House CRANBOURNE'
$ws.Cells.Item(19, 1).Value = 'li: Elwood Children''s Centre ELWOOD'
$ws.Cells.Item(20, 1).Value = 'li: Fernwood Fitness Clayton CLAYTON'
$ws.Cells.Item(21, 1).Value = 'li: Fernwood Fitness Narre Warren NARRE WARREN'
$ws.Cells.Item(22, 1).Value = 'li: GESAC Childcare BENTLEIGH EAST'
$ws.Cells.Item(23, 1).Value = 'li: Glen Iris Creche GLEN IRIS'
$ws.Cells.Item(24, 1).Value = 'li: Goodlife Chelsea Heights CHELSEA HEIGHTS'
$ws.Cells.Item(25, 1).Value = 'li: Goodlife Fountain Gate NARRE WARREN'
$ws.Cells.Item(26, 1).Value = 'li: Goodlife Karingal KARINGAL'
$ws.Cells.Item(27, 1).Value = 'li: Hallam Community Centre Inc HALLAM'
$ws.Cells.Item(28, 1).Value = 'li: Hope Frankston Heights FRANKSTON'
$ws.Cells.Item(29, 1).Value = 'li: Melbourne Grammar School, Grimwade House CAULFIELD'
$ws.Cells.Item(30, 1).Value = 'li: Merkaz Bentleigh BENTLEIGH'
$ws.Cells.Item(31, 1).Value = 'li: Niño Early Learning Adventures - Malvern East MALVERN EAST'
$ws.Cells.Item(32, 1).Value = 'li: Only About Children Cheltenham CHELTENHAM'
$ws.Cells.Item(33, 1).Value = 'li: Overport Primary School OSHC - Extend FRANKSTON'
$ws.Cells.Item(34, 1).Value = 'li: Paisley Park Early Learning Centre Cranbourne CRANBOURNE'
$ws.Cells.Item(35, 1).Value = 'li: Pavillion - Frankston & District Netball Association Inc FRANKSTON'
$ws.Cells.Item(36, 1).Value = 'li: Peninsula Aquatic Recreation Centre FRANKSTON'
$ws.Cells.Item(37, 1).Value = 'li: Play Zone - Highett HIGHETT'
$ws.Cells.Item(38, 1).Value = 'li: Playtime Warragul WARRAGUL'
$ws.Cells.Item(39, 1).Value = 'li: Sacre Coeur Oshclub GLEN IRIS'
$ws.Cells.Item(40, 1).Value = 'li: Sandybeach Centre SANDRINGHAM'
$ws.Cells.Item(41, 1).Value = 'li: Smaland Springvale Ikea SPRINGVALE'
$ws.Cells.Item(42, 1).Value = 'li: Somerville Recreation Centre Child Care SOMERVILLE'
$ws.Cells.Item(43, 1).Value = 'li: Toorak Primary School OSHClub TOORAK'
$ws.Cells.Item(44, 1).Value = 'li: Upper Beaconsfield Community Early Learning Centre BEACONSFIELD UPPER'
$ws.Cells.Item(45, 1).Value = 'li: Yavneh College TheirCare ELSTERNWICK'
$ws.Cells.Item(46, 1).Value = 'Schools closedTheDepartment hasbeen advised of the followingschool closures:'
$ws.Cells.Item(47, 1).Value = 'li: Dandenong High School, DANDENONG'
$ws.Cells.Item(48, 1).Value = 'li: Gleneagles Secondary College, ENDEAVOUR HILLS'
$ws.Cells.Item(49, 1).Value = 'li: Leibler Yavneh College, ELSTERNWICK'
$ws.Cells.Item(50, 1).Value = 'li: Overport Primary School, FRANKSTON'
$ws.Cells.Item(51, 1).Value = 'TAFEThe Department hasnotbeen advised of any TAFE closures.North-Eastern Victoria RegionEarly childhood servicesThe Department has been advised of the following early childhood service closures:'
$ws.Cells.Item(52, 1).Value = 'li: Aqualink Box Hill Creche BOX HILL'
$ws.Cells.Item(53, 1).Value = 'li: Aqualink Nunawading Creche FOREST HILL'
$ws.Cells.Item(54, 1).Value = 'li: Camp Australia - Camberwell Boys Grammar Junior School OSHC CANTERBURY'
$ws.Cells.Item(55, 1).Value = 'li: Camp Australia - Mountain Gate Primary School OSHC FERNTREE GULLY'
$ws.Cells.Item(56, 1).Value = 'li: Camp Australia - St Clement of Rome School OSHC BULLEEN'
$ws.Cells.Item(57, 1).Value = 'li: Camp Australia - Strathcona Baptist Girls Junior School OSHC CANTERBURY'
$ws.Cells.Item(58, 1).Value = 'li: Camp Australia - Templestowe Heights Primary School OSHC TEMPLESTOWE LOWER'
$ws.Cells.Item(59, 1).Value = 'li: Clever Kids Childcare - Ashburton ASHBURTON'
$ws.Cells.Item(60, 1).Value = 'li: Fitness First Doncaster (Playzone) DONCASTER'
$ws.Cells.Item(61, 1).Value = 'li: MakerDojo HAWTHORN'
$ws.Cells.Item(62, 1).Value = 'li: Roberts McCubbin OSHClub BOX HILL SOUTH'
$ws.Cells.Item(63, 1).Value = 'li: St Andrews Christian College Outside School Hours Care WANTIRNA SOUTH'
$ws.Cells.Item(64, 1).Value = 'li: Victoria Road OSHClub LILYDALE'
$ws.Cells.Item(65, 1).Value = 'li: Wattle Park Children''s Services Centre BURWOOD'
$ws.Cells.Item(66, 1).Value = 'li: West Hawthorn Early Childhood Centre HAWTHORN'
$ws.Cells.Item(67, 1).Value = 'Schools closedThe Department hasbeen advisedof the followingschool closures:'
$ws.Cells.Item(68, 1).Value = 'li: Canterbury Girl''s Secondary College, CANTERBURY'
$ws.Cells.Item(69, 1).Value = 'li: Glen Iris Primary School, GLEN IRIS'
$ws.Cells.Item(70, 1).Value = 'li: Kew High School, KEW'
$ws.Cells.Item(71, 1).Value = 'li: Roberts McCubbin Primary School, BOX HILL SOUTH'
$ws.Cells.Item(72, 1).Value = 'TAFEThe Department hasnotbeen advised of any TAFE closures.North-Western Victoria RegionEarly childhood servicesThe Department has been advised of the following early childhood service closures:'
$ws.Cells.Item(73, 1).Value = 'li: Bright Stars Early Years Child Care Centre EPPING'
$ws.Cells.Item(74, 1).Value = 'li: Camp Australia - Ivanhoe Grammar OSHC IVANHOE'
$ws.Cells.Item(75, 1).Value = 'li: Camp Australia - Ivanhoe Grammar Plenty Valley Campus OSHC MERNDA'
$ws.Cells.Item(76, 1).Value = 'li: Collingwood College Afterschool Care and Vacation Care Program COLLINGWOOD'
$ws.Cells.Item(77, 1).Value = 'li: Diamond Creek Community Centre DIAMOND CREEK'
$ws.Cells.Item(78, 1).Value = 'li: Diamond Valley Sports and Fitness Centre GREENSBOROUGH'
$ws.Cells.Item(79, 1).Value = 'li: Eltham Leisure Centre ELTHAM'
$ws.Cells.Item(80, 1).Value = 'li: Gumboots Early Learning Sunbury SUNBURY'
$ws.Cells.Item(81, 1).Value = 'li: Kangaroo Ground Primary Combined OSHC KANGAROO GROUND'
$ws.Cells.Item(82, 1).Value = 'li: New Futures Broadmeadows BROADMEADOWS'
$ws.Cells.Item(83, 1).Value = 'li: New Futures Epping EPPING'
$ws.Cells.Item(84, 1).Value = 'li: Nino Early Learning Adventures - Bundoora BUNDOORA'
$ws.Cells.Item(85, 1).Value = 'li: Only About Children Coonans Hill PASCOE VALE SOUTH'
$ws.Cells.Item(86, 1).Value = 'li: Richmond West Afterschool Care and Vacation Care Program RICHMOND'
$ws.Cells.Item(87, 1).Value = 'SchoolsThe Department hasbeen advised of the following school closures:'
$ws.Cells.Item(88, 1).Value = 'li: Al Siraat College, EPPING'
$ws.Cells.Item(89, 1).Value = 'li: Charles La Trobe, MACLEOD WEST'
$ws.Cells.Item(90, 1).Value = 'li: Diamond Valley College, DIAMOND REEK'
$ws.Cells.Item(91, 1).Value = 'li: Epping Secondary College, EPPING'
$ws.Cells.Item(92, 1).Value = 'li: Fitzroy High School, FITZROY'
$ws.Cells.Item(93, 1).Value = 'li: Fitzroy Primary School, FITZROY'
$ws.Cells.Item(94, 1).Value = 'li: Gladstone Park Secondary College'
$ws.Cells.Item(95, 1).Value = 'li: Macleod College, MACLEOD'
$ws.Cells.Item(96, 1).Value = 'li: Moomba Primary School, FAWKNER'
$ws.Cells.Item(97, 1).Value = 'li: Newbury Primary School, CRAGIEBURN'
$ws.Cells.Item(98, 1).Value = 'li: Northcote High School, NORTHCOTE'
$ws.Cells.Item(99, 1).Value = 'li: Northern College of Arts and Technology (NCAT), PRESTON'
$ws.Cells.Item(100, 1).Value = 'li: Pascoe Vale Girls Secondary College, OAK PARK'
$ws.Cells.Item(101, 1).Value = 'li: Penders Grove Primary School, THORNBURY'
$ws.Cells.Item(102, 1).Value = 'li: Princes Hill Secondary College, PRINCES HILL'
$ws.Cells.Item(103, 1).Value = 'li: Roxburgh College, ROXBURGH PARK'
$ws.Cells.Item(104, 1).Value = 'TAFEThe Department hasnotbeen advised of any TAFE closures.South-Western Victoria RegionThe Department has not been advised of any school, early childhood service or TAFE closures, or buscancellations.Early childhood servicesThe Department has been advised of the following early childhood service closures:'
$ws.Cells.Item(105, 1).Value = 'li: Aerotots Activity Centre WERRIBEE'
$ws.Cells.Item(106, 1).Value = 'li: Altona Meadows Community Centre Occasional Care ALTONA MEADOWS'
$ws.Cells.Item(107, 1).Value = 'li: Amici - Westbourne Early Learning Centre HOPPERS CROSSING'
$ws.Cells.Item(108, 1).Value = 'li: Aquapulse Creche HOPPERS CROSSING'
$ws.Cells.Item(109, 1).Value = 'li: Big Childcare - Manor Lakes P-12 College OSHC WYNDHAM VALE'
$ws.Cells.Item(110, 1).Value = 'li: Blackwood Street Neighbourhood House YARRAVILLE'
$ws.Cells.Item(111, 1).Value = 'li: Bluewater Leisure Centre Creche COLAC'
$ws.Cells.Item(112, 1).Value = 'li: Camp Australia - Baden Powell P-9 College Derrimut Heath Campus OSHC HOPPERS CROSSING'
$ws.Cells.Item(113, 1).Value = 'li: Camp Australia - Footscray City Primary School OSHC FOOTSCRAY'
$ws.Cells.Item(114, 1).Value = 'li: Camp Australia - Haileybury City Campus OSHC WEST MELBOURNE'
$ws.Cells.Item(115, 1).Value = 'li: Camp Australia - Melton Christian College OSHC Melton South'
$ws.Cells.Item(116, 1).Value = 'li: Carranballac Jamieson OSHClub POINT COOK'
$ws.Cells.Item(117, 1).Value = 'li: Eagle Stadium WERRIBEE'
$ws.Cells.Item(118, 1).Value = 'li: Energy Force Fitness Creche DRYSDALE'
$ws.Cells.Item(119, 1).Value = 'li: Explorers Early Learning - Williams Landing WILLIAMS LANDING'
$ws.Cells.Item(120, 1).Value = 'li: Fernwood Fitness Sydenham SYDENHAM'
$ws.Cells.Item(121, 1).Value = 'li: Footscray PS TheirCare FOOTSCRAY'
$ws.Cells.Item(122, 1).Value = 'li: Footscray West PS TheirCare WEST FOOTSCRAY'
$ws.Cells.Item(123, 1).Value = 'li: Fun 4 All Occasional Care Center WERRIBEE'
$ws.Cells.Item(124, 1).Value = 'li: Genesis Maidstone MAIDSTONE'
$ws.Cells.Item(125, 1).Value = 'li: Glen Gala Children''s Centre SUNSHINE WEST'
$ws.Cells.Item(126, 1).Value = 'li: Goodlife Essendon Child Minding ESSENDON'
$ws.Cells.Item(127, 1).Value = 'li: Goodlife Geelong BELMONT'
$ws.Cells.Item(128, 1).Value = 'li: Goodlife Point Cook POINT COOK'
$ws.Cells.Item(129, 1).Value = 'li: Goodlife Taylors Lakes TAYLORS LAKES'
$ws.Cells.Item(130, 1).Value = 'li: Goodstart Early Learning Grovedale - Pioneer Road GROVEDALE'
$ws.Cells.Item(131, 1).Value = 'li: Goodstart Early Learning Kealba KEALBA'
$ws.Cells.Item(132, 1).Value = 'li: Hopetoun Early Years Centre FLEMINGTON'
$ws.Cells.Item(133, 1).Value = 'li: Keilor Basketball Stadium Creche KEILOR PARK'
$ws.Cells.Item(134, 1).Value = 'li: Kensington Neighbourhood House Inc KENSINGTON'
$ws.Cells.Item(135, 1).Value = 'li: Kids Club Kensington Early Learning Centre KENSINGTON'
$ws.Cells.Item(136, 1).Value = 'li: Maribyrnong Aquatic Centre Occasional Child Care MARIBYRNONG'
$ws.Cells.Item(137, 1).Value = 'li: New Futures Braybrook BRAYBROOK'
$ws.Cells.Item(138, 1).Value = 'li: Shuter Street Occasional Care MOONEE PONDS'
$ws.Cells.Item(139, 1).Value = 'li: South Kingsville Community Centre SOUTH KINGSVILLE'
$ws.Cells.Item(140, 1).Value = 'li: St Vincent de Paul TheirCare STRATHMORE'
$ws.Cells.Item(141, 1).Value = 'li: Sunshine Leisure Centre SUNSHINE'
$ws.Cells.Item(142, 1).Value = 'li: Wallaby Childcare Harpley Werribee'
$ws.Cells.Item(143, 1).Value = 'li: Westgate Sports & Leisure Centre ALTONA NORTH'
$ws.Cells.Item(144, 1).Value = 'li: Willaura Primary School OSHC ARARAT'
$ws.Cells.Item(145, 1).Value = 'li: Williamstown Community and Education Centre Occasional Care WILLIAMSTOWN'
$ws.Cells.Item(146, 1).Value = 'li: Winifred Nance Kindergarten COLAC'
$ws.Cells.Item(147, 1).Value = 'li: Yara Childcare Centres TRUGANINA'
$ws.Cells.Item(148, 1).Value = 'li: Yarraville Community Centre YARRAVILLE'
$ws.Cells.Item(149, 1).Value = 'Schools closedThe Department has been advised of the following school closures:'
$ws.Cells.Item(150, 1).Value = 'li: Al Taqwa College, TRUGANINA'
$ws.Cells.Item(151, 1).Value = 'li: Catholic Regional College, SYDENHAM'
$ws.Cells.Item(152, 1).Value = 'li: Clonard Secondary College, GEELONG WEST'
$ws.Cells.Item(153, 1).Value = 'li: Copperfield College, DELAHEY'
$ws.Cells.Item(154, 1).Value = 'li: Footscray Primary School and OSHC, FOOTSCRAY'
$ws.Cells.Item(155, 1).Value = 'li: Footscray West Primary School, FOOTSCARY WEST'
$ws.Cells.Item(156, 1).Value = 'li: Grovedale West Primary School, GROVEDALE'
$ws.Cells.Item(157, 1).Value = 'li: Lowther Hall Anglican Grammar School, ESSENDON'
$ws.Cells.Item(158, 1).Value = 'li: Manor Lakes P-12 College, WYNDHAM VALE'
$ws.Cells.Item(159, 1).Value = 'li: Parkville College (Parkville and Malmsbury Campus)'
$ws.Cells.Item(160, 1).Value = 'li: St Brigid''s College, HORSHAM'
$ws.Cells.Item(161, 1).Value = 'li: Sydenham Hillside Primary School, SYDENHAM'
$ws.Cells.Item(162, 1).Value = 'li: Tarneit Senior College, TARNEIT'
$ws.Cells.Item(163, 1).Value = 'li: Thomas Carr College, TARNEIT'
$ws.Cells.Item(164, 1).Value = 'li: Trinity College, COLAC'
$ws.Cells.Item(165, 1).Value = 'li: Victorian College of the Arts, ALBERT PARK'
$ws.Cells.Item(166, 1).Value = 'TAFEThe Department hasnotbeen advised of any TAFE closures.The Department hasnotbeen advised of any school, early childhood service or TAFE closures, or buscancellations.Bus service cancellations or alterationsFor Term 2 2020, schools bus services will continue to be provided to support student travel to schools where needed.With learning from home arrangements, local principals are authorised to cancel or suspend school buses where not required.Families and students are encouraged to contact their school directly to enquire if their bus service is continuing to be provided.Find more about bus services:School bus services in Term 2 (for schools)School bus services in Term 2 (for parents)VicRoadsRegional roadsBack to topLast Update: 23 July 2020Website navigationFor parentsEarly childhoodSchoolsTAFE and trainingThe Department<li>' + $cr + '                                 <a href="/about/educationstate">Education State' + $cr + '                                 ' + $cr + '                                 </a>' + $cr + '                             </li>Support linksHelp in your languageAccessibilityPrivacyCopyrightContact linksContact usReport a website issue'
$ws.Cells.Item(167, 1).Value = 'li: State Government of Victoria, Australia © 2019'
$ws.Cells.Item(168, 1).Value = 'We respectfully acknowledge the Traditional Owners of country throughout Victoria and pay respect to the ongoing living cultures of First Peoples.Our website uses a free tool to translate into other languages. This tool is a guide and may not be accurate. For more, see:Information in your languageOld Search Code PlaceholdCludo Code for www.education.vic.gov.auHotjar Tracking Code for www.education.vic.gov.auSR-1181393'
Write-Host "Done setting all rows."
